# Append the new quarterly-data row (2021-Q3 / "01-07-2021") to Sheet1,
# mirroring the "Actualización desde MV -datos-" update: one new row 76
# with a date-label in column A and 40 numeric observations in B:AP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 76

# Column A holds a text label that looks like a date ("01-07-2021").
# Assigning it straight to .Value would make Excel auto-convert it into a
# real date serial. Instead, enter it as a text formula and immediately
# paste-special just the value back onto itself; this keeps the stored
# cell a plain shared string (matching the existing A2:A75 cells) without
# leaving any residual number-format/style behind on the cell.
$ws.Cells.Item($row, 1).Formula = '="01-07-2021"'
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)  # xlPasteValues

$values = @(23194,3585,1482,1393,139,377,75,55,62,506,493,12,10011,8600,792,410,151,58,138,105,33,110,993,52,601,123,217,7049,6161,313,128,99,44,304,409,163,246,160,58,101,234)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # B=2 ... AP=42
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
